# Update automatico via Actualizar 02-04-2021 22-10-22
#
# 1) Refresh the "last checked" timestamp for the previous batch of rows
#    (156:169) to the new recalculated value.
# 2) Append a fresh batch of 14 availability rows (170:183) with the new
#    run's timestamp, mirroring the existing Nombre/URL/Disponibilidad/Fecha
#    layout and the column-B hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = 44231.90265305556

# --- 1) refresh D156:D169 -------------------------------------------------
for ($r = 156; $r -le 169; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldDate
}

# --- 2) append rows 170:183 -----------------------------------------------
$newDate = 44231.92382483089

$data = @(
    @{ Row=170; A="Odoo";               B="https://www.dataintelligence-group.com/";                     C="Disponible" },
    @{ Row=171; A="Blackbox";           B="https://serviciodashboard.azurewebsites.net/";                C="Disponible" },
    @{ Row=172; A="PowerBI";            B="https://powerbi.microsoft.com/es-es/";                        C="Disponible" },
    @{ Row=173; A="Dropbox";            B="https://www.dropbox.com/";                                    C="Disponible" },
    @{ Row=174; A="Odoo";               B="https://dataintelligence.store/";                             C="Disponible" },
    @{ Row=175; A="GEE";                B="https://app-data-i.users.earthengine.app/";                   C="Disponible" },
    @{ Row=176; A="UtilidadesOdoo";     B="https://odooutil.azurewebsites.net/";                         C="Disponible" },
    @{ Row=177; A="Filtros Dashboard";  B="https://filtradordashboard.azurewebsites.net/";                C="Disponible" },
    @{ Row=178; A="MapStore";           B="https://ide.dataintelligence-group.com/mapstore/#/";           C="Disponible" },
    @{ Row=179; A="GeoServer";          B="https://ide.dataintelligence-group.com/geoserver/web/?0";      C="Disponible" },
    @{ Row=180; A="Tomcat";             B="https://ide.dataintelligence-group.com/";                      C="Disponible" },
    @{ Row=181; A="Shiny";              B="https://rpubs.com/dataintelligence/";                          C="Disponible" },
    @{ Row=182; A="Github";             B="https://github.com/Sud-Austral/";                              C="Disponible" },
    @{ Row=183; A="EZ Exporter";        B="https://ezexporter.highviewapps.com/exports/export-profile/";  C="Disponible" }
)

foreach ($item in $data) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 3).Value = $item.C

    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $item.B

    $hashIdx = $item.B.IndexOf("#")
    if ($hashIdx -ge 0) {
        $address = $item.B.Substring(0, $hashIdx)
        $subAddress = $item.B.Substring($hashIdx + 1)
        $ws.Hyperlinks.Add($bCell, $address, $subAddress) | Out-Null
    } else {
        $ws.Hyperlinks.Add($bCell, $item.B) | Out-Null
    }

    # Hyperlinks.Add() stamps its own style xf; reapply the sheet's usual
    # "Hyperlink" cell style afterwards so it matches the existing B2:B169
    # cells (style index 2) instead of a freshly minted one.
    $bCell.Style = "Hyperlink"
}
